# Update num.txt | 2025-12-10 08:02:24 JST
# Adds the newest committee meeting (第634回 / 2025-12-10) as a new row directly
# below the header row, pushing all existing meeting rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row right after the header row (row 1).
# Every pre-existing meeting row (formerly rows 2-15) shifts down to rows 3-16.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the 第634回 meeting information.
$ws.Cells.Item(2, 1).Value = "第634回"
$ws.Cells.Item(2, 2).Value = "2025年12月10日（令和7年12月10日）"
$ws.Cells.Item(2, 3).Value = "１費用対効果評価専門組織からの報告について`n２令和７年度補正予算案の閣議決定について`n３令和８年度診療報酬改定に関する基本的な見解（各号意見）について`n"
$ws.Cells.Item(2, 4).Value = "－"
$ws.Cells.Item(2, 5).Value = "資料`n"
$ws.Cells.Item(2, 6).Value = "－"
